$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '330.21'
Set-TextValue 'E2' '1.25%'
Set-TextValue 'G2' '11'
Set-TextValue 'D3' '44.37'
Set-TextValue 'E3' '0.24%'
Set-TextValue 'G3' '11'
Set-TextValue 'D4' '5.488'
Set-TextValue 'E4' '-1.64%'
Set-TextValue 'G4' '11'
Set-TextValue 'D5' '0.08023'
Set-TextValue 'E5' '-0.26%'
Set-TextValue 'G5' '11'
Set-TextValue 'D6' '2.115'
Set-TextValue 'E6' '12.01%'
Set-TextValue 'G6' '11'
Set-TextValue 'B7' 'BTSEToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D7' '2.586'
Set-TextValue 'E7' '-2.11%'
Set-TextValue 'G7' '11'
Set-TextValue 'B8' 'MXToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D8' '0.9525'
Set-TextValue 'E8' '0.98%'
Set-TextValue 'G8' '11'
Set-TextValue 'D9' '0.1144'
Set-TextValue 'E9' '-1.30%'
Set-TextValue 'G9' '11'
Set-TextValue 'D10' '0.1896'
Set-TextValue 'E10' '3.02%'
Set-TextValue 'G10' '11'
Set-TextValue 'D11' '10.72'
Set-TextValue 'E11' '28.32%'
Set-TextValue 'G11' '11'
Set-TextValue 'D12' '0.09876'
Set-TextValue 'E12' '1.04%'
Set-TextValue 'G12' '11'
Set-TextValue 'D13' '0.04818'
Set-TextValue 'E13' '10.23%'
Set-TextValue 'G13' '11'
Set-TextValue 'D14' '0.1064'
Set-TextValue 'E14' '-0.11%'
Set-TextValue 'G14' '11'
Set-TextValue 'D15' '0.001267'
Set-TextValue 'E15' '-1.04%'
Set-TextValue 'G15' '11'
Set-TextValue 'E16' '-3.33%'
Set-TextValue 'G16' '11'
Set-TextValue 'D17' '0.005935'
Set-TextValue 'E17' '-0.18%'
Set-TextValue 'G17' '11'
Set-TextValue 'E18' '-6.65%'
Set-TextValue 'G18' '11'
Set-TextValue 'D19' '4.404'
Set-TextValue 'E19' '2.52%'
Set-TextValue 'G19' '11'
Set-TextValue 'D20' '0.3459'
Set-TextValue 'E20' '-1.03%'
Set-TextValue 'G20' '11'
Set-TextValue 'D21' '0.1399'
Set-TextValue 'E21' '1.48%'
Set-TextValue 'G21' '11'
Set-TextValue 'D22' '0.2503'
Set-TextValue 'E22' '-5.59%'
Set-TextValue 'G22' '11'
Set-TextValue 'E23' '2.41%'
Set-TextValue 'G23' '11'
Set-TextValue 'D24' '0.004370'
Set-TextValue 'E24' '-3.04%'
Set-TextValue 'G24' '11'
Set-TextValue 'E25' '-4.86%'
Set-TextValue 'G25' '11'
Set-TextValue 'D26' '0.0003743'
Set-TextValue 'E26' '-6.25%'
Set-TextValue 'G26' '11'
Set-TextValue 'G27' '11'
Set-TextValue 'G28' '11'
Set-TextValue 'G29' '11'
Set-TextValue 'G30' '11'
Set-TextValue 'G31' '11'
Set-TextValue 'G32' '11'
Set-TextValue 'G33' '11'
Set-TextValue 'G34' '11'
Set-TextValue 'G35' '11'
Set-TextValue 'G36' '11'
Set-TextValue 'G37' '11'
Set-TextValue 'D38' '0.02596'
Set-TextValue 'E38' '-0.87%'
Set-TextValue 'G38' '11'
Set-TextValue 'D39' '0.05816'
Set-TextValue 'E39' '6.68%'
Set-TextValue 'G39' '11'
Set-TextValue 'D40' '0.007564'
Set-TextValue 'E40' '-0.35%'
Set-TextValue 'G40' '11'
Set-TextValue 'D41' '0.1404'
Set-TextValue 'E41' '0.68%'
Set-TextValue 'G41' '11'
Set-TextValue 'D42' '0.007146'
Set-TextValue 'E42' '-2.61%'
Set-TextValue 'G42' '11'
Set-TextValue 'E43' '-0.15%'
Set-TextValue 'G43' '11'
Set-TextValue 'D44' '0.008801'
Set-TextValue 'E44' '-0.30%'
Set-TextValue 'G44' '11'
Set-TextValue 'D45' '0.00007006'
Set-TextValue 'E45' '1.21%'
Set-TextValue 'G45' '11'
Set-TextValue 'E46' '-0.10%'
Set-TextValue 'G46' '11'
Set-TextValue 'D47' '0.0005798'
Set-TextValue 'E47' '-0.23%'
Set-TextValue 'G47' '11'
Set-TextValue 'D48' '0.003529'
Set-TextValue 'E48' '55.33%'
Set-TextValue 'G48' '11'
Set-TextValue 'D49' '0.003498'
Set-TextValue 'E49' '-4.90%'
Set-TextValue 'G49' '11'
Set-TextValue 'E50' '-0.10%'
Set-TextValue 'G50' '11'
Set-TextValue 'E51' '-0.10%'
Set-TextValue 'G51' '11'
